# Updates cryptos list values (Price column D, Volume(1h) column E)
# for Sheet1, per the data refresh described in the commit message.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "68.042.50"
$ws.Cells.Item(2, 5).Value = "  +1.41%  "
$ws.Cells.Item(3, 4).Value = "3.263.71"
$ws.Cells.Item(3, 5).Value = "  -0.11%  "
$ws.Cells.Item(5, 4).Value = "'586.08"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +1.08%  "
$ws.Cells.Item(6, 4).Value = "'184.64"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +4.08%  "
$ws.Cells.Item(7, 5).Value = "  -0.05%  "
$ws.Cells.Item(8, 5).Value = "  -1.01%  "
$ws.Cells.Item(9, 5).Value = "  +3.59%  "
$ws.Cells.Item(10, 4).Value = "'6.70"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  -0.75%  "
$ws.Cells.Item(11, 5).Value = "  +0.05%  "
$ws.Cells.Item(12, 4).Value = "3.829.33"
$ws.Cells.Item(12, 5).Value = "  -0.24%  "
$ws.Cells.Item(13, 5).Value = "  +0.27%  "
$ws.Cells.Item(14, 5).Value = "  +1.55%  "
$ws.Cells.Item(15, 4).Value = "68.068.24"
$ws.Cells.Item(15, 5).Value = "  +1.40%  "
$ws.Cells.Item(16, 5).Value = "  +2.58%  "
$ws.Cells.Item(17, 4).Value = "3.260.99"
$ws.Cells.Item(17, 5).Value = "  -0.25%  "
$ws.Cells.Item(18, 5).Value = "  -0.18%  "
$ws.Cells.Item(20, 4).Value = "'381.83"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +2.61%  "
$ws.Cells.Item(21, 4).Value = "'7.69"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  +0.59%  "
$ws.Cells.Item(22, 5).Value = "  +0.04%  "
$ws.Cells.Item(23, 4).Value = "'71.38"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +0.11%  "
$ws.Cells.Item(25, 4).Value = "'0.0000121"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +1.01%  "
$ws.Cells.Item(26, 4).Value = "'9.90"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  +0.50%  "
$ws.Cells.Item(27, 5).Value = "  +2.53%  "
$ws.Cells.Item(28, 5).Value = "  +0.05%  "
$ws.Cells.Item(29, 5).Value = "  +0.18%  "
$ws.Cells.Item(30, 5).Value = "  +0.64%  "
$ws.Cells.Item(31, 4).Value = "'7.28"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  +6.53%  "
$ws.Cells.Item(32, 4).Value = "'22.90"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +1.08%  "
$ws.Cells.Item(33, 5).Value = "  +0.04%  "
$ws.Cells.Item(34, 5).Value = "  +0.33%  "
$ws.Cells.Item(35, 5).Value = "  +2.43%  "
$ws.Cells.Item(36, 4).Value = "'162.70"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  -3.28%  "
$ws.Cells.Item(37, 5).Value = "  -0.34%  "
$ws.Cells.Item(38, 4).Value = "'0.836"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  -2.59%  "
$ws.Cells.Item(39, 5).Value = "  +5.23%  "
$ws.Cells.Item(40, 4).Value = "'26.61"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -2.66%  "
$ws.Cells.Item(42, 4).Value = "'2.60"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  +0.48%  "
$ws.Cells.Item(43, 4).Value = "'41.39"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  +2.09%  "
$ws.Cells.Item(44, 4).Value = "'25.52"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +2.05%  "
$ws.Cells.Item(45, 4).Value = "'347.40"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -0.36%  "
$ws.Cells.Item(46, 4).Value = "2.652.74"
$ws.Cells.Item(46, 5).Value = "  -4.03%  "
$ws.Cells.Item(47, 5).Value = "  +1.28%  "
$ws.Cells.Item(48, 4).Value = "'0.0285"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +1.62%  "
$ws.Cells.Item(49, 4).Value = "'32.08"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +5.21%  "
$ws.Cells.Item(50, 5).Value = "  -1.18%  "
$ws.Cells.Item(51, 4).Value = "'0.998"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +1.06%  "
